# Update of regular expressions for identifying tracking cookies
# -> The "Mac Address" row (row 6) gets its example value filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# B6 (next to "Mac Address" in A6) gets the example MAC address.
# This adds a new shared string and keeps the existing yellow-fill style (s="2").
$ws.Range("B6").Value = "00:00:0c:07:ac:0e"

# Move the selection/view down to B8, scrolled so row 4 is the first visible row.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B8").Select()
